# Applies the NATMI TPM re-run update to Ndp-Fzd4.xlsx:
#  - adds a new 'ECs' sending-cluster block (rows 2-5) and shifts the
#    original MuSCs block down (rows 6-9)
#  - adds a new 'Resolving-Mac' target-cluster row to each block
#  - refreshes all numeric LR-score columns (E:T) with the new TPM values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
  @{ row=2; A="ECs"; B="Ndp"; C="Fzd4"; D="ECs"; E=1; F=0.3333333333333333; G=0.008202333333333334; H=0.024607; I=0.0618542966452333; J=0.08999806156895879; K=3; L=1; M=19.524618; N=58.573854; O=0.4154885426712971; P=0.4539723485554654; Q=0.160147425042; R=1.441326825378; S=0.02569975157108609; T=0.04085663137589959 },
  @{ row=3; A="ECs"; B="Ndp"; C="Fzd4"; D="FAPs"; E=1; F=0.3333333333333333; G=0.008202333333333334; H=0.024607; I=0.0618542966452333; J=0.08999806156895879; K=3; L=1; M=15.24435933333334; N=45.73307800000001; O=0.324403614112412; P=0.3544508583357054; Q=0.1250393167051111; R=1.125353850346; S=0.02006575738009492; T=0.0318998901716671 },
  @{ row=4; A="ECs"; B="Ndp"; C="Fzd4"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.008202333333333334; H=0.024607; I=0.0618542966452333; J=0.08999806156895879; K=2; L=1; M=11.9507005; N=23.901401; O=0.2543137660693869; P=0.1852460510065796; Q=0.09802362906783334; R=0.588141774407; S=0.01573039912742233; T=0.01667178550389663 },
  @{ row=5; A="ECs"; B="Ndp"; C="Fzd4"; D="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.008202333333333334; H=0.024607; I=0.0618542966452333; J=0.08999806156895879; K=2; L=0.6666666666666666; M=0.272275; N=0.816825; O=0.005794077146903843; P=0.006330742102249548; Q=0.002233290308333333; R=0.020099612775; S=0.0003583885666299573; T=0.0005697545174954544 },
  @{ row=6; A="MuSCs"; B="Ndp"; C="Fzd4"; D="ECs"; E=2; F=1; G=0.124405; H=0.24881; I=0.9381457033547667; J=0.9100019384310413; K=3; L=1; M=19.524618; N=58.573854; O=0.4154885426712971; P=0.4539723485554654; Q=2.42896010229; R=14.57376061374; S=0.389788791100211; T=0.4131157171795659 },
  @{ row=7; A="MuSCs"; B="Ndp"; C="Fzd4"; D="FAPs"; E=2; F=1; G=0.124405; H=0.24881; I=0.9381457033547667; J=0.9100019384310413; K=3; L=1; M=15.24435933333334; N=45.73307800000001; O=0.324403614112412; P=0.3544508583357054; Q=1.896474522863334; R=11.37884713718; S=0.3043378567323171; T=0.3225509681640383 },
  @{ row=8; A="MuSCs"; B="Ndp"; C="Fzd4"; D="MuSCs"; E=2; F=1; G=0.124405; H=0.24881; I=0.9381457033547667; J=0.9100019384310413; K=2; L=1; M=11.9507005; N=23.901401; O=0.2543137660693869; P=0.1852460510065796; Q=1.4867268957025; R=5.94690758281; S=0.2385833669419646; T=0.168574265502683 },
  @{ row=9; A="MuSCs"; B="Ndp"; C="Fzd4"; D="Resolving-Mac"; E=2; F=1; G=0.124405; H=0.24881; I=0.9381457033547667; J=0.9100019384310413; K=2; L=0.6666666666666666; M=0.272275; N=0.816825; O=0.005794077146903843; P=0.006330742102249548; Q=0.03387237137499999; R=0.20323422825; S=0.005435688580273886; T=0.005760987584754094 }
)

foreach ($r in $rowData) {
  $ws.Range("A" + $r.row).Value = $r.A
  $ws.Range("B" + $r.row).Value = $r.B
  $ws.Range("C" + $r.row).Value = $r.C
  $ws.Range("D" + $r.row).Value = $r.D
  $ws.Range("E" + $r.row).Value = $r.E
  $ws.Range("F" + $r.row).Value = $r.F
  $ws.Range("G" + $r.row).Value = $r.G
  $ws.Range("H" + $r.row).Value = $r.H
  $ws.Range("I" + $r.row).Value = $r.I
  $ws.Range("J" + $r.row).Value = $r.J
  $ws.Range("K" + $r.row).Value = $r.K
  $ws.Range("L" + $r.row).Value = $r.L
  $ws.Range("M" + $r.row).Value = $r.M
  $ws.Range("N" + $r.row).Value = $r.N
  $ws.Range("O" + $r.row).Value = $r.O
  $ws.Range("P" + $r.row).Value = $r.P
  $ws.Range("Q" + $r.row).Value = $r.Q
  $ws.Range("R" + $r.row).Value = $r.R
  $ws.Range("S" + $r.row).Value = $r.S
  $ws.Range("T" + $r.row).Value = $r.T
}
